$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Cells.Item(2, 4)
$c.NumberFormat = "@"
$c.Value = "51.586.83"
$c.Style = "Normal"
$ws.Cells.Item(2, 5).Value = "  +1.21%  "

$c = $ws.Cells.Item(3, 4)
$c.NumberFormat = "@"
$c.Value = "3.016.44"
$c.Style = "Normal"
$ws.Cells.Item(3, 5).Value = "  +2.18%  "

$ws.Cells.Item(4, 5).Value = "  +0.14%  "

$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = "@"
$c.Value = "379.56"
$c.Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  -0.17%  "

$c = $ws.Cells.Item(6, 4)
$c.NumberFormat = "@"
$c.Value = "102.34"
$c.Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  +0.49%  "

$c = $ws.Cells.Item(7, 4)
$c.NumberFormat = "@"
$c.Value = "0.545"
$c.Style = "Normal"
$ws.Cells.Item(7, 5).Value = "  +0.79%  "

$ws.Cells.Item(8, 5).Value = "  +0.05%  "

$c = $ws.Cells.Item(9, 4)
$c.NumberFormat = "@"
$c.Value = "0.589"
$c.Style = "Normal"
$ws.Cells.Item(9, 5).Value = "  +0.71%  "

$c = $ws.Cells.Item(10, 4)
$c.NumberFormat = "@"
$c.Value = "36.66"
$c.Style = "Normal"
$ws.Cells.Item(10, 5).Value = "  +1.15%  "

$ws.Cells.Item(11, 5).Value = "  -0.17%  "

$ws.Cells.Item(12, 5).Value = "  +1.33%  "

$c = $ws.Cells.Item(13, 4)
$c.NumberFormat = "@"
$c.Value = "3.491.78"
$c.Style = "Normal"
$ws.Cells.Item(13, 5).Value = "  +2.43%  "

$c = $ws.Cells.Item(14, 4)
$c.NumberFormat = "@"
$c.Value = "18.41"
$c.Style = "Normal"
$ws.Cells.Item(14, 5).Value = "  +0.07%  "

$c = $ws.Cells.Item(15, 4)
$c.NumberFormat = "@"
$c.Value = "7.69"
$c.Style = "Normal"
$ws.Cells.Item(15, 5).Value = "  -0.78%  "

$c = $ws.Cells.Item(16, 4)
$c.NumberFormat = "@"
$c.Value = "3.008.57"
$c.Style = "Normal"
$ws.Cells.Item(16, 5).Value = "  +2.03%  "

$c = $ws.Cells.Item(17, 4)
$c.NumberFormat = "@"
$c.Value = "0.973"
$c.Style = "Normal"
$ws.Cells.Item(17, 5).Value = "  -3.55%  "

$c = $ws.Cells.Item(18, 4)
$c.NumberFormat = "@"
$c.Value = "10.53"
$c.Style = "Normal"
$ws.Cells.Item(18, 5).Value = "  -14.72%  "

$c = $ws.Cells.Item(19, 4)
$c.NumberFormat = "@"
$c.Value = "51.575.98"
$c.Style = "Normal"
$ws.Cells.Item(19, 5).Value = "  +1.16%  "

$ws.Cells.Item(20, 5).Value = "  -0.06%  "

$ws.Cells.Item(21, 5).Value = "  +0.49%  "

$c = $ws.Cells.Item(22, 4)
$c.NumberFormat = "@"
$c.Value = "0.0₃0960"
$c.Style = "Normal"
$ws.Cells.Item(22, 5).Value = "  +0.76%  "

$c = $ws.Cells.Item(23, 4)
$c.NumberFormat = "@"
$c.Value = "69.89"
$c.Style = "Normal"
$ws.Cells.Item(23, 5).Value = "  +0.47%  "

$c = $ws.Cells.Item(24, 4)
$c.NumberFormat = "@"
$c.Value = "265.78"
$c.Style = "Normal"
$ws.Cells.Item(24, 5).Value = "  -0.25%  "

$ws.Cells.Item(25, 5).Value = "  -7.21%  "

$c = $ws.Cells.Item(26, 4)
$c.NumberFormat = "@"
$c.Value = "8.21"
$c.Style = "Normal"
$ws.Cells.Item(26, 5).Value = "  +2.60%  "

$c = $ws.Cells.Item(27, 4)
$c.NumberFormat = "@"
$c.Value = "7.60"
$c.Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  +10.02%  "

$c = $ws.Cells.Item(28, 4)
$c.NumberFormat = "@"
$c.Value = "0.172"
$c.Style = "Normal"
$ws.Cells.Item(28, 5).Value = "  +4.85%  "

$ws.Cells.Item(29, 5).Value = "  +0.12%  "

$c = $ws.Cells.Item(30, 4)
$c.NumberFormat = "@"
$c.Value = "26.14"
$c.Style = "Normal"
$ws.Cells.Item(30, 5).Value = "  +1.46%  "

$c = $ws.Cells.Item(31, 4)
$c.NumberFormat = "@"
$c.Value = "0.108"
$c.Style = "Normal"
$ws.Cells.Item(31, 5).Value = "  +0.89%  "

$c = $ws.Cells.Item(32, 4)
$c.NumberFormat = "@"
$c.Value = "10.23"
$c.Style = "Normal"
$ws.Cells.Item(32, 5).Value = "  -2.08%  "

$ws.Cells.Item(33, 5).Value = "  +0.29%  "

$c = $ws.Cells.Item(34, 4)
$c.NumberFormat = "@"
$c.Value = "50.46"
$c.Style = "Normal"
$ws.Cells.Item(34, 5).Value = "  -0.46%  "

$c = $ws.Cells.Item(35, 4)
$c.NumberFormat = "@"
$c.Value = "33.66"
$c.Style = "Normal"
$ws.Cells.Item(35, 5).Value = "  -0.69%  "

$c = $ws.Cells.Item(36, 4)
$c.NumberFormat = "@"
$c.Value = "0.0447"
$c.Style = "Normal"
$ws.Cells.Item(36, 5).Value = "  +3.43%  "

$ws.Cells.Item(37, 5).Value = "  -0.10%  "

$ws.Cells.Item(38, 5).Value = "  +3.34%  "

$c = $ws.Cells.Item(39, 4)
$c.NumberFormat = "@"
$c.Value = "0.297"
$c.Style = "Normal"
$ws.Cells.Item(39, 5).Value = "  +14.91%  "

$c = $ws.Cells.Item(40, 4)
$c.NumberFormat = "@"
$c.Value = "16.91"
$c.Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  +1.74%  "

$c = $ws.Cells.Item(41, 4)
$c.NumberFormat = "@"
$c.Value = "1.85"
$c.Style = "Normal"
$ws.Cells.Item(41, 5).Value = "  +1.87%  "

$c = $ws.Cells.Item(42, 4)
$c.NumberFormat = "@"
$c.Value = "127.68"
$c.Style = "Normal"
$ws.Cells.Item(42, 5).Value = "  +6.80%  "

$ws.Cells.Item(43, 5).Value = "  -0.72%  "

$ws.Cells.Item(44, 5).Value = "  +2.83%  "

$ws.Cells.Item(45, 5).Value = "  +5.87%  "

$c = $ws.Cells.Item(46, 4)
$c.NumberFormat = "@"
$c.Value = "21.60"
$c.Style = "Normal"
$ws.Cells.Item(46, 5).Value = "  +1.22%  "

$c = $ws.Cells.Item(47, 4)
$c.NumberFormat = "@"
$c.Value = "2.09"
$c.Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  +3.56%  "

$ws.Cells.Item(48, 5).Value = "  +2.74%  "

$c = $ws.Cells.Item(49, 4)
$c.NumberFormat = "@"
$c.Value = "2.025.75"
$c.Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  +0.13%  "

$c = $ws.Cells.Item(50, 4)
$c.NumberFormat = "@"
$c.Value = "3.319.72"
$c.Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  +2.36%  "

$c = $ws.Cells.Item(51, 4)
$c.NumberFormat = "@"
$c.Value = "0.0319"
$c.Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  +0.58%  "
